$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Selling, General & Administrative
$ws.Range("E6").Value = "'-256"
$ws.Range("F6").Value = "'-264"
$ws.Range("G6").Value = "'-272"

# Row 8: Total Operating Expenses
$ws.Range("E8").Value = "'-422"
$ws.Range("F8").Value = "'-439"
$ws.Range("G8").Value = "'-469"

# Row 9: EBITDA
$ws.Range("E9").Value = "'1,559"
$ws.Range("F9").Value = "'1,832"
$ws.Range("G9").Value = "'2,143"

# Row 10: Operating Income (EBIT)
$ws.Range("E10").Value = "'1,523"
$ws.Range("F10").Value = "'1,795"
$ws.Range("G10").Value = "'2,105"

# Row 12: Earnings Before Tax
$ws.Range("E12").Value = "'1,460"
$ws.Range("F12").Value = "'1,728"
$ws.Range("G12").Value = "'2,036"

# Row 14: Net Income
$ws.Range("E14").Value = "'1,352"
$ws.Range("F14").Value = "'1,599"
$ws.Range("G14").Value = "'1,880"
